# Applies the "lab.cell.*" (Cell/battery) translation-key additions plus a couple
# of follow-up keys (lab.mod.table.power, lab.mod.preview.power, lab.liquid.table.vendor)
# to the "Import" sheet of the translations workbook.
#
# Each new row follows the existing "cs | <key> | <translation>" layout already
# used for every other row on this sheet (column A is always the literal "cs").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$firstNewRow = 1081

$rows = @(
    @('lab.cell.label', 'Články'),
    @('lab.cell.title', 'Články'),
    @('lab.cell.button.create', 'Nový článek'),
    @('lab.cell.filter.title', 'Filtrování článků'),
    @('lab.cell.table.name', 'Jméno'),
    @('lab.cell.table.size', 'Velikost'),
    @('lab.cell.table.drain', 'Vybíjecí proud'),
    @('lab.cell.table.voltage', 'Napětí'),
    @('lab.cell.table.ohm', 'Bezpečný odpor'),
    @('lab.cell.table.vendor', 'Výrobce'),
    @('lab.cell.table.footer.label', 'Počet článků [{{data.total}}]'),
    @('lab.cell.name.label', 'Jméno'),
    @('lab.cell.name.label.tooltip', 'Použijte prosím obchodní označení článku tak, aby ostatní mohli tento článek snadno najít. '),
    @('lab.cell.drain.label', 'Vybíjecí proud'),
    @('lab.cell.drain.label.tooltip', 'Bezpečný vybíjecí proud článku; uveďte prosím reálnou hodnotu, poněvadž aplikace s ní pak bude počítat pro různé kontroly a pokud by tato hodnota neodpovádala realitě, může vás to pak ohrozit. Tak či tak, pro výpočet bezpečného odporu na tomto článkui je použito pouze 75% vybíjecího proudu, poněvadž se očekává, že výrobci budou optimisti a nám to za to nestojí.'),
    @('lab.cell.vendorId.label', 'Výrobce'),
    @('lab.cell.voltage.label', 'Napětí'),
    @('lab.cell.voltage.label.tooltip', 'Tuto hodnotu standardně není třeba měnit, pokud ovšem nemáte nějaký ultra zvláštní článek.'),
    @('lab.cell.create.submit', 'Uložit'),
    @('lab.cell.size.label', 'Velikost článku'),
    @('lab.cell.size.label.tooltip', 'Použijte prosím standardní rozměr udaný výrobcem (např. 21700) .'),
    @('lab.cell.size.label.required', 'Prosím udejte velikost článku.'),
    @('lab.cell.create.success', 'Článek [{{data.vendor.name}} {{data.name}}] byl úspěšně uložen.'),
    @('lab.cell.preview.name', 'Jméno'),
    @('lab.cell.context.menu', 'Článek [{{data.vendor.name}} {{data.name}}]'),
    @('lab.cell.preview', 'Náhled článku'),
    @('lab.cell.button.edit', 'Upravit článek'),
    @('lab.cell.button.delete', 'Odstranit článek'),
    @('lab.cell.button.delete.confirm.title', 'Odstranit článek'),
    @('lab.cell.button.delete.confirm.ok', 'Odstranit článek'),
    @('lab.cell.button.delete.confirm', 'Opravdu si přejete odstranit vybraný článek? Tato skce je nezvratná a může být smazáno větší množství dat.'),
    @('lab.cell.deleted.success', 'Článek [{{data.vendor.name}} {{data.name}}] byl úspěšně odstraněn.'),
    @('lab.cell.index.label', 'Článek [{{data.vendor.name}} {{data.name}}]'),
    @('lab.cell.index.title', 'Detail článku'),
    @('lab.cell.preview.size', 'Velikost'),
    @('lab.cell.preview.drain', 'Vybíjecí proud'),
    @('lab.cell.preview.ohm', 'Minimální bezpečný odpor'),
    @('lab.cell.preview.voltage', 'Pracovní napětí'),
    @('lab.cell.update.submit', 'Aktualizovat'),
    @('lab.cell.updated.message', 'Článek [{{data.vendor.name}} {{data.name}}] byl úspěšně aktualizován.'),
    @('lab.mod.table.power', 'Výkon'),
    @('lab.mod.preview.power', 'Maximální výkon'),
    @('lab.liquid.table.vendor', 'Výrobce')
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $firstNewRow + $i
    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $rows[$i][0]
    $ws.Cells.Item($r, 3).Value = $rows[$i][1]
}

$lastNewRow = $firstNewRow + $rows.Count - 1

# Match the formatting already used by every other data row (wrapped text,
# 10pt font) so the new rows render/print the same as the rest of the sheet.
$newRange = $ws.Range("A" + $firstNewRow + ":C" + $lastNewRow)
$newRange.WrapText = $true
$newRange.Font.Size = 10

# Leave the selection where the author left it after adding the new rows
# (partway through typing the new block, not on the very last row).
$ws.Range("B1112").Select() | Out-Null
